$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.308.95"
$ws.Range("E2").Value = "  -2.93%  "

$ws.Range("D3").Value = "1.831.73"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.77%  "

$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5159"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3215"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06725"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7638"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07674"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").Value = "1.811.78"
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.009"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.77%  "

$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007878"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.88%  "

$ws.Range("D20").Value = "26.356.49"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").Value = "2.078.15"
$ws.Range("E21").Value = "  -2.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.398"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.896"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.21%  "

$ws.Range("E25").Value = "  -3.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.648"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("E28").Value = "  -3.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.172"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08702"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04822"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.845"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6782"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.091"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01775"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.184"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4876"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.45%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.27%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8924"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.662"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4173"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1252"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.065"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05878"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.62%  "
